$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 56, shifting rows
# 56..68 down to 57..69 (matches the new dimension A1:R69).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(56, 1).Value = 3
$ws.Cells.Item(56, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(56, 3).Value = "Coquimbo"
$ws.Cells.Item(56, 4).Value = 44782
$ws.Cells.Item(56, 5).Value = 5
$ws.Cells.Item(56, 6).Value = 100112035
$ws.Cells.Item(56, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 55
$ws.Cells.Item(56, 11).Value = 15000
$ws.Cells.Item(56, 12).Value = 15000
$ws.Cells.Item(56, 13).Value = 15000
$ws.Cells.Item(56, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(56, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(56, 16).Value = 1000
$ws.Cells.Item(56, 17).Value = 15
$ws.Cells.Item(56, 18).Value = "Hortaliza"
